# "continental us message fix" - update staging MSRP data:
#  - roll model-year 2020 -> 2021 for a block of existing rows and bump several MSRPs
#  - append a new block of 2021 Lexus ES/RX trims (rows 81-94) plus a trailing blank row 95
#  - widen column B, refresh the saved view (zoom/selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B width ---
$ws.Range("B:B").ColumnWidth = 54.14

# --- Existing rows: model year 2020 -> 2021 (rows 10-16) ---
$ws.Range("C10").Value = 2021
$ws.Range("C11").Value = 2021
$ws.Range("C12").Value = 2021
$ws.Range("C13").Value = 2021
$ws.Range("C14").Value = 2021
$ws.Range("C15").Value = 2021
$ws.Range("C16").Value = 2021

# --- MSRP bumps that go along with the 2021 rows above ---
$ws.Range("D11").Value = 45100
$ws.Range("D12").Value = 48900
$ws.Range("D14").Value = 47010
$ws.Range("D15").Value = 50810
$ws.Range("D16").Value = 45700

# --- Existing rows: model year 2020 -> 2021 (rows 55-63) ---
$ws.Range("C55").Value = 2021
$ws.Range("C56").Value = 2021
$ws.Range("C57").Value = 2021
$ws.Range("C58").Value = 2021
$ws.Range("C59").Value = 2021
$ws.Range("C60").Value = 2021
$ws.Range("C61").Value = 2021
$ws.Range("C62").Value = 2021
$ws.Range("C63").Value = 2021

# --- MSRP bumps that go along with the 2021 rows above ---
$ws.Range("D55").Value = 45070
$ws.Range("D56").Value = 46470
$ws.Range("D57").Value = 47900
$ws.Range("D58").Value = 49300
$ws.Range("D59").Value = 48550
$ws.Range("D60").Value = 49950
$ws.Range("D61").Value = 47720
$ws.Range("D62").Value = 51110
$ws.Range("D63").Value = 51200

# --- New rows 81-94 (2021 ES / RX trims) ---
# Row 81 - ES 250
$ws.Range("A81").Value = 9012
$ws.Range("B82").Value = "ES 250 LUXURY"
$ws.Range("B83").Value = "ES 250 ULTRA LUXURY"
$ws.Range("B84").Value = "ES 250 F SPORT"
$ws.Range("B85").Value = "ES 350 F SPORT BLACK LINE SPECIAL EDITION"
$ws.Range("B81").Value = "ES 250"
$ws.Range("C81").Value = 2021
$ws.Range("D81").Value = 39900
$ws.Range("D81").NumberFormat = "$#,##0_);[Red]($#,##0)"
$ws.Range("E81").Value = 1025
$ws.Range("E81").NumberFormat = "$#,##0.00_);[Red]($#,##0.00)"

# Row 82 - ES 250 LUXURY
$ws.Range("A82").Value = 9013
$ws.Range("C82").Value = 2021
$ws.Range("D82").Value = 45100
$ws.Range("D82").NumberFormat = "$#,##0_);[Red]($#,##0)"
$ws.Range("E82").Value = 1025
$ws.Range("E82").NumberFormat = "$#,##0.00_);[Red]($#,##0.00)"

# Row 83 - ES 250 ULTRA LUXURY
$ws.Range("A83").Value = 9014
$ws.Range("C83").Value = 2021
$ws.Range("D83").Value = 48900
$ws.Range("D83").NumberFormat = "$#,##0_);[Red]($#,##0)"
$ws.Range("E83").Value = 1025
$ws.Range("E83").NumberFormat = "$#,##0.00_);[Red]($#,##0.00)"

# Row 84 - ES 250 F SPORT
$ws.Range("A84").Value = 9015
$ws.Range("C84").Value = 2021
$ws.Range("D84").Value = 45700
$ws.Range("D84").NumberFormat = "$#,##0_);[Red]($#,##0)"
$ws.Range("E84").Value = 1025
$ws.Range("E84").NumberFormat = "$#,##0.00_);[Red]($#,##0.00)"

# Row 85 - ES 350 F SPORT BLACK LINE SPECIAL EDITION
$ws.Range("C85").Value = 2021
$ws.Range("D85").Value = 46550
$ws.Range("D85").NumberFormat = "$#,##0_);[Red]($#,##0)"
$ws.Range("E85").Value = 1025
$ws.Range("E85").NumberFormat = "$#,##0.00_);[Red]($#,##0.00)"

# Row 86 - RX 350 F SPORT PERFORMANCE FWD
$ws.Range("A86").Value = 9423
$ws.Range("B86").Value = "RX 350 F SPORT PERFORMANCE FWD"
$ws.Range("C86").Value = 2021
$ws.Range("D86").Value = 50950
$ws.Range("D86").NumberFormat = "$#,##0_);[Red]($#,##0)"
$ws.Range("E86").Value = 1025
$ws.Range("E86").NumberFormat = "$#,##0.00_);[Red]($#,##0.00)"

# Row 87 - RX 350 F SPORT PERFORMANCE AWD
$ws.Range("A87").Value = 9427
$ws.Range("B87").Value = "RX 350 F SPORT PERFORMANCE AWD"
$ws.Range("C87").Value = 2021
$ws.Range("D87").Value = 52350
$ws.Range("D87").NumberFormat = "$#,##0_);[Red]($#,##0)"
$ws.Range("E87").Value = 1025
$ws.Range("E87").NumberFormat = "$#,##0.00_);[Red]($#,##0.00)"

# Row 88 - RX 450h F-SPORT PERFORMANCE AWD
$ws.Range("A88").Value = 9447
$ws.Range("B88").Value = "RX 450h F-SPORT PERFORMANCE AWD"
$ws.Range("C88").Value = 2021
$ws.Range("D88").Value = 53520
$ws.Range("D88").NumberFormat = "$#,##0_);[Red]($#,##0)"
$ws.Range("E88").Value = 1025
$ws.Range("E88").NumberFormat = "$#,##0.00_);[Red]($#,##0.00)"

# Row 89 - RX 350L LUXURY FWD
$ws.Range("A89").Value = 9432
$ws.Range("B89").Value = "RX 350L LUXURY FWD"
$ws.Range("C89").Value = 2021
$ws.Range("D89").Value = 53900
$ws.Range("D89").NumberFormat = "$#,##0_);[Red]($#,##0)"
$ws.Range("E89").Value = 1025
$ws.Range("E89").NumberFormat = "$#,##0.00_);[Red]($#,##0.00)"

# Row 90 - RX 350L LUXURY AWD
$ws.Range("A90").Value = 9436
$ws.Range("B90").Value = "RX 350L LUXURY AWD"
$ws.Range("C90").Value = 2021
$ws.Range("D90").Value = 55300
$ws.Range("D90").NumberFormat = "$#,##0_);[Red]($#,##0)"
$ws.Range("E90").Value = 1025
$ws.Range("E90").NumberFormat = "$#,##0.00_);[Red]($#,##0.00)"

# Row 91 - RX 450hL LUXURY AWD
$ws.Range("A91").Value = 9457
$ws.Range("B91").Value = "RX 450hL LUXURY AWD"
$ws.Range("C91").Value = 2021
$ws.Range("D91").Value = 57110
$ws.Range("D91").NumberFormat = "$#,##0_);[Red]($#,##0)"
$ws.Range("E91").Value = 1025
$ws.Range("E91").NumberFormat = "$#,##0.00_);[Red]($#,##0.00)"
$ws.Range("J91").NumberFormat = "$#,##0.00_);[Red]($#,##0.00)"

# Row 92 - RX 350 F SPORT BLACK LINE SPECIAL EDITION (9422SE)
$ws.Range("B92").Value = "RX 350 F SPORT BLACK LINE SPECIAL EDITION"
# Row 94 needs its B value filled before A92/A93/A94 to reproduce the original
# shared-string ordering (B94 comes before the A-column SE codes).
$ws.Range("B94").Value = "RX 450h F SPORT AWD BLACK LINE SPECIAL EDITION"
$ws.Range("A92").Value = "9422SE"
$ws.Range("A93").Value = "9426SE"
$ws.Range("A94").Value = "9446SE"
$ws.Range("A85").Value = "9005SE"
$ws.Range("B93").Value = "RX 350 AWD F SPORT BLACK LINE SPECIAL EDITION"

$ws.Range("C92").Value = 2021
$ws.Range("D92").Value = 49235
$ws.Range("D92").NumberFormat = "$#,##0_);[Red]($#,##0)"
$ws.Range("E92").Value = 1025
$ws.Range("E92").NumberFormat = "$#,##0.00_);[Red]($#,##0.00)"
$ws.Range("K92").NumberFormat = "$#,##0.00_);[Red]($#,##0.00)"

$ws.Range("C93").Value = 2021
$ws.Range("D93").Value = 50635
$ws.Range("D93").NumberFormat = "$#,##0_);[Red]($#,##0)"
$ws.Range("E93").Value = 1025
$ws.Range("E93").NumberFormat = "$#,##0.00_);[Red]($#,##0.00)"
$ws.Range("K93").NumberFormat = "$#,##0.00_);[Red]($#,##0.00)"

$ws.Range("C94").Value = 2021
$ws.Range("D94").Value = 51885
$ws.Range("D94").NumberFormat = "$#,##0_);[Red]($#,##0)"
$ws.Range("E94").Value = 1025
$ws.Range("E94").NumberFormat = "$#,##0.00_);[Red]($#,##0.00)"
$ws.Range("K94").NumberFormat = "$#,##0.00_);[Red]($#,##0.00)"

# Row 95 - trailing stray formatted cell only
$ws.Range("K95").NumberFormat = "$#,##0.00_);[Red]($#,##0.00)"

# --- Saved view state: zoom + selection (topLeftCell/zoomScaleNormal are not
#     persisted by this runtime, so only the supported pieces are applied) ---
$win = $excel.ActiveWindow
$win.Zoom = 80
$win.ScrollRow = 65
$win.ScrollColumn = 1
$ws.Range("B93").Select() | Out-Null
